$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OriginTemplate")

# Update the F2 header cell: "SecretARN" -> "SecretID"
$ws.Range("F2").Value = "SecretID"

# Update the big instructions cell (A1): rename the "SecretARN" bullet to "SecretID"
$description = @"
该模版文件用于批量创建SDP的数据源，请不要修改表头(彩色部分)的信息和顺序。具体规范如下如下：
- InstanceName：必填，不能重复。填写数据库实例名称。
- SSL: 必填。1表示启用SSL，0表示不启用。
- Description: 选填。
- JDBC_URL: 必填。数据库实例 格式：jdbc:protocol://host:port。或者，数据库实例中的数据库(database) 格式：jdbc:protocol://host:port/database
- JDBC_Databases: 选填。如果JDBC_URL填的是mysql协议，此项可保留为空。如果，JDBC_URL填写的是其他数据库实例，请填写相同实例下的数据库中待检测的database，多个database用半角逗号隔开。
- SecretID: 选填，secret的ARN。您可以指定SDP所在的AWS账号下的secret manager服务中的secret Id (请提前将用户名/密码存储在此secret中)。
- Username: 必填。（如果secretID填了，那么此参数将被忽略）
- Password: 必填。（如果secretID填了，那么此参数将被忽略）
- AccountID: 请提前在SDP平台添加好账户（account ID），如果账号不存在将报错。
- Region: 请提前在SDP平台添加好区域（region），如果区域不存在将报错。
- ProviderID：请提前在SDP平台添加好供应商(Provider)，如果供应商不存在将报错。
---------------------------------------------------------------------------------------------------------------------------------------------------------------------------------------------------
上传模版后，等待批量处理完，您可以看到具体的批处理报告。
- 错误的情况：假设10条数据源，其中8条信息正确，2条不正确。那么，上传模版后，SDP平台将成功生成8个数据源，同时报错。您可以下载报告文件。报告文件中将提示您哪2条有错误，并附带错误说明。您可以对错误的2条数据源信息进行修改，并再次提交。
- 数据源已经存在的情况：如果模版中定义的数据源，与SDP平台上已经存在，那么，上传后会提示报错。如果您依旧希望添加模版中的数据源，请您先在SDP界面上删除对应的数据源，再次上传。
---------------------------------------------------------------------------------------------------------------------------------------------------------------------------------------------------
更多的说明，详见文档。https://awslabs.github.io/sensitive-data-protection-on-aws/zh/user-guide/data-catalog-create-jdbc/
"@
$ws.Range("A1").Value = $description
